{"js": "// Replace the division-problem answer text in each table cell with the\n// newly generated values, per the diff. Each \"before\" string is unique in\n// the document, so a targeted search+replace per pair is safe and keeps\n// run formatting (font, size, etc.) intact.\nconst replacements = [\n  [\"750\u00f75=150, 0\", \"636\u00f74=159, 0\"],\n  [\"847\u00f72=423, 1\", \"697\u00f74=174, 1\"],\n  [\"626\u00f75=125, 1\", \"487\u00f75=97, 2\"],\n  [\"340\u00f77=48, 4\", \"767\u00f74=191, 3\"],\n  [\"915\u00f72=457, 1\", \"976\u00f78=122, 0\"],\n  [\"770\u00f74=192, 2\", \"881\u00f78=110, 1\"],\n  [\"941\u00f76=156, 5\", \"491\u00f77=70, 1\"],\n  [\"105\u00f77=15, 0\", \"669\u00f76=111, 3\"],\n  [\"337\u00f79=37, 4\", \"300\u00f78=37, 4\"],\n  [\"862\u00f76=143, 4\", \"580\u00f73=193, 1\"],\n  [\"443\u00f78=55, 3\", \"631\u00f78=78, 7\"],\n  [\"207\u00f77=29, 4\", \"771\u00f75=154, 1\"],\n  [\"472\u00f75=94, 2\", \"985\u00f77=140, 5\"],\n  [\"372\u00f76=62, 0\", \"574\u00f77=82, 0\"],\n  [\"484\u00f73=161, 1\", \"438\u00f75=87, 3\"],\n  [\"534\u00f79=59, 3\", \"572\u00f74=143, 0\"],\n  [\"687\u00f79=76, 3\", \"753\u00f72=376, 1\"],\n  [\"877\u00f73=292, 1\", \"463\u00f79=51, 4\"],\n  [\"423\u00f72=211, 1\", \"897\u00f72=448, 1\"],\n  [\"874\u00f78=109, 2\", \"951\u00f77=135, 6\"],\n  [\"568\u00f77=81, 1\", \"761\u00f78=95, 1\"],\n  [\"244\u00f74=61, 0\", \"338\u00f78=42, 2\"],\n  [\"633\u00f77=90, 3\", \"322\u00f79=35, 7\"],\n  [\"580\u00f79=64, 4\", \"799\u00f73=266, 1\"],\n  [\"544\u00f72=272, 0\", \"565\u00f73=188, 1\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the division-problem answer text in each table cell with the\n# newly generated values, per the diff. Each \"before\" string is unique in\n# the document, so Find/Replace per pair is safe and keeps run formatting\n# (font, size, etc.) intact since Find.Execute only swaps the text inside\n# the existing run(s).\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"750\u00f75=150, 0\", \"636\u00f74=159, 0\"),\n    @(\"847\u00f72=423, 1\", \"697\u00f74=174, 1\"),\n    @(\"626\u00f75=125, 1\", \"487\u00f75=97, 2\"),\n    @(\"340\u00f77=48, 4\", \"767\u00f74=191, 3\"),\n    @(\"915\u00f72=457, 1\", \"976\u00f78=122, 0\"),\n    @(\"770\u00f74=192, 2\", \"881\u00f78=110, 1\"),\n    @(\"941\u00f76=156, 5\", \"491\u00f77=70, 1\"),\n    @(\"105\u00f77=15, 0\", \"669\u00f76=111, 3\"),\n    @(\"337\u00f79=37, 4\", \"300\u00f78=37, 4\"),\n    @(\"862\u00f76=143, 4\", \"580\u00f73=193, 1\"),\n    @(\"443\u00f78=55, 3\", \"631\u00f78=78, 7\"),\n    @(\"207\u00f77=29, 4\", \"771\u00f75=154, 1\"),\n    @(\"472\u00f75=94, 2\", \"985\u00f77=140, 5\"),\n    @(\"372\u00f76=62, 0\", \"574\u00f77=82, 0\"),\n    @(\"484\u00f73=161, 1\", \"438\u00f75=87, 3\"),\n    @(\"534\u00f79=59, 3\", \"572\u00f74=143, 0\"),\n    @(\"687\u00f79=76, 3\", \"753\u00f72=376, 1\"),\n    @(\"877\u00f73=292, 1\", \"463\u00f79=51, 4\"),\n    @(\"423\u00f72=211, 1\", \"897\u00f72=448, 1\"),\n    @(\"874\u00f78=109, 2\", \"951\u00f77=135, 6\"),\n    @(\"568\u00f77=81, 1\", \"761\u00f78=95, 1\"),\n    @(\"244\u00f74=61, 0\", \"338\u00f78=42, 2\"),\n    @(\"633\u00f77=90, 3\", \"322\u00f79=35, 7\"),\n    @(\"580\u00f79=64, 4\", \"799\u00f73=266, 1\"),\n    @(\"544\u00f72=272, 0\", \"565\u00f73=188, 1\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
